$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Serial")
$tbl = $ws1.ListObjects.Add(1, $ws1.Range("B2:C33"), $null, 1)
Write-Host $tbl.Name
$tbl.TableStyle = "TableStyleLight9"
